# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (value
# "stock" for the existing data row), inserted right after the "total"
# column and before the "date" column. This pushes the existing
# date / legislator_name / legislator_id columns one position to the
# right (H->I, I->J, J->K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H, shifting former H:J (date, legislator_name,
# legislator_id) to I:K.
$ws.Columns("H:H").Insert()

# New column header + value.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
